$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (Miami @ Boston) entirely - the consensus tracker now has a single game row
$ws.Rows(3).Delete()

# Update remaining row 2 (Denver @ Miami) with refreshed consensus figures
$ws.Range("D2").Value = 219
$ws.Range("E2").Value = -8.5
$ws.Range("F2").Value = "Denver"
$ws.Range("G2").Value = "Miami"
$ws.Range("H2").Value = 0.5520833333333334
$ws.Range("I2").Value = 0.4444444444444444
$ws.Range("J2").Value = 115.7926829268293
$ws.Range("K2").Value = 109.4756097560976
$ws.Range("L2").Value = 97.49512195121947
$ws.Range("M2").Value = 95.51707317073166
$ws.Range("N2").Value = 118.3817073170732
$ws.Range("O2").Value = 113.9914634146342
$ws.Range("P2").Value = 115.0475609756098
$ws.Range("Q2").Value = 114.4048780487804
$ws.Range("R2").Value = 76.57439024390247
$ws.Range("S2").Value = 77.84146341463416
$ws.Range("T2").Value = 0.3610853658536586
$ws.Range("U2").Value = 0.4089268292682925
$ws.Range("V2").Value = 0.6026829268292681
$ws.Range("W2").Value = 0.5751097560975612
$ws.Range("X2").Value = 0.2604024390243902
$ws.Range("Y2").Value = 0.2742073170731708
$ws.Range("Z2").Value = 12.5280487804878
$ws.Range("AA2").Value = 11.81463414634147
$ws.Range("AB2").Value = 11.67195121951219
$ws.Range("AC2").Value = 13.95121951219512
$ws.Range("AD2").Value = 0.1996524390243903
$ws.Range("AE2").Value = 0.214859756097561
$ws.Range("AF2").Value = 1.009526442256576
$ws.Range("AG2").Value = 0.9544516979607459
$ws.Range("AH2").Value = 0.9730033350886431
$ws.Range("AI2").Value = 1.096134566113401
$ws.Range("AJ2").Value = 11.16695670584409
$ws.Range("AK2").Value = 10.58096036232879
$ws.Range("AL2").Value = 0.6463414634146342
$ws.Range("AM2").Value = 0.5365853658536586
$ws.Range("AN2").Value = 51.5
$ws.Range("AO2").Value = 49.5
$ws.Range("AP2").Value = 75.09999999999999
$ws.Range("AQ2").Value = 76.2
$ws.Range("AR2").Value = 0.68
$ws.Range("AS2").Value = 0.6153846153846154
$ws.Range("AT2").Value = 0.4922378056669106
$ws.Range("AU2").Value = 0.503512542778542
$ws.Range("AV2").Value = 2.16
$ws.Range("AW2").Value = 4.23
$ws.Range("AX2").Value = 0
$ws.Range("AY2").Value = 0.1352191751241632
$ws.Range("AZ2").Value = 0
$ws.Range("BA2").Value = 0.1187381436314363
